$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb2"
$ws.Cells.Item(2, 3).Value = "Ephb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 41.519606
$ws.Cells.Item(2, 8).Value = 124.558818
$ws.Cells.Item(2, 9).Value = 0.7305114279806179
$ws.Cells.Item(2, 10).Value = 0.7630546295388222
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.875567333333333
$ws.Cells.Item(2, 14).Value = 5.626702
$ws.Cells.Item(2, 15).Value = 0.9147184316015459
$ws.Cells.Item(2, 16).Value = 0.9397527619538806
$ws.Cells.Item(2, 17).Value = 77.87281670647067
$ws.Cells.Item(2, 18).Value = 700.855350358236
$ws.Cells.Item(2, 19).Value = 0.6682122676694364
$ws.Cells.Item(2, 20).Value = 0.7170826956308033

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb2"
$ws.Cells.Item(3, 3).Value = "Ephb1"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 41.519606
$ws.Cells.Item(3, 8).Value = 124.558818
$ws.Cells.Item(3, 9).Value = 0.7305114279806179
$ws.Cells.Item(3, 10).Value = 0.7630546295388222
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.010998
$ws.Cells.Item(3, 14).Value = 0.032994
$ws.Cells.Item(3, 15).Value = 0.005363749481003509
$ws.Cells.Item(3, 16).Value = 0.005510546431623061
$ws.Cells.Item(3, 17).Value = 0.4566326267880001
$ws.Cells.Item(3, 18).Value = 4.109693641092
$ws.Cells.Item(3, 19).Value = 0.003918280292698171
$ws.Cells.Item(3, 20).Value = 0.004204847965938613

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb2"
$ws.Cells.Item(4, 3).Value = "Ephb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 41.519606
$ws.Cells.Item(4, 8).Value = 124.558818
$ws.Cells.Item(4, 9).Value = 0.7305114279806179
$ws.Cells.Item(4, 10).Value = 0.7630546295388222
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.163866
$ws.Cells.Item(4, 14).Value = 0.327732
$ws.Cells.Item(4, 15).Value = 0.07991781891745053
$ws.Cells.Item(4, 16).Value = 0.05473669161449624
$ws.Cells.Item(4, 17).Value = 6.803651756796001
$ws.Cells.Item(4, 18).Value = 40.821910540776
$ws.Cells.Item(4, 19).Value = 0.05838088001848322
$ws.Cells.Item(4, 20).Value = 0.04176708594208019

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efnb2"
$ws.Cells.Item(5, 3).Value = "Ephb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 7.183879666666667
$ws.Cells.Item(5, 8).Value = 21.551639
$ws.Cells.Item(5, 9).Value = 0.1263958572665066
$ws.Cells.Item(5, 10).Value = 0.1320266053993819
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.875567333333333
$ws.Cells.Item(5, 14).Value = 5.626702
$ws.Cells.Item(5, 15).Value = 0.9147184316015459
$ws.Cells.Item(5, 16).Value = 0.9397527619538806
$ws.Cells.Item(5, 17).Value = 13.47385002939756
$ws.Cells.Item(5, 18).Value = 121.264650264578
$ws.Cells.Item(5, 19).Value = 0.1156166203197518
$ws.Cells.Item(5, 20).Value = 0.1240723670754643

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb2"
$ws.Cells.Item(6, 3).Value = "Ephb1"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7.183879666666667
$ws.Cells.Item(6, 8).Value = 21.551639
$ws.Cells.Item(6, 9).Value = 0.1263958572665066
$ws.Cells.Item(6, 10).Value = 0.1320266053993819
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.010998
$ws.Cells.Item(6, 14).Value = 0.032994
$ws.Cells.Item(6, 15).Value = 0.005363749481003509
$ws.Cells.Item(6, 16).Value = 0.005510546431623061
$ws.Cells.Item(6, 17).Value = 0.07900830857400001
$ws.Cells.Item(6, 18).Value = 0.7110747771660001
$ws.Cells.Item(6, 19).Value = 0.0006779557138142185
$ws.Cells.Item(6, 20).Value = 0.0007275387392628701

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb2"
$ws.Cells.Item(7, 3).Value = "Ephb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7.183879666666667
$ws.Cells.Item(7, 8).Value = 21.551639
$ws.Cells.Item(7, 9).Value = 0.1263958572665066
$ws.Cells.Item(7, 10).Value = 0.1320266053993819
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.163866
$ws.Cells.Item(7, 14).Value = 0.327732
$ws.Cells.Item(7, 15).Value = 0.07991781891745053
$ws.Cells.Item(7, 16).Value = 0.05473669161449624
$ws.Cells.Item(7, 17).Value = 1.177193625458
$ws.Cells.Item(7, 18).Value = 7.063161752748001
$ws.Cells.Item(7, 19).Value = 0.0101012812329406
$ws.Cells.Item(7, 20).Value = 0.007226699584654754

$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Efnb2"
$ws.Cells.Item(8, 3).Value = "Ephb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.229894
$ws.Cells.Item(8, 8).Value = 0.6896819999999999
$ws.Cells.Item(8, 9).Value = 0.004044840748830231
$ws.Cells.Item(8, 10).Value = 0.004225032410066655
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.875567333333333
$ws.Cells.Item(8, 14).Value = 5.626702
$ws.Cells.Item(8, 15).Value = 0.9147184316015459
$ws.Cells.Item(8, 16).Value = 0.9397527619538806
$ws.Cells.Item(8, 17).Value = 0.4311816765293333
$ws.Cells.Item(8, 18).Value = 3.880635088764
$ws.Cells.Item(8, 19).Value = 0.003699890385848011
$ws.Cells.Item(8, 20).Value = 0.0039704858767048

$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Efnb2"
$ws.Cells.Item(9, 3).Value = "Ephb1"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.229894
$ws.Cells.Item(9, 8).Value = 0.6896819999999999
$ws.Cells.Item(9, 9).Value = 0.004044840748830231
$ws.Cells.Item(9, 10).Value = 0.004225032410066655
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.010998
$ws.Cells.Item(9, 14).Value = 0.032994
$ws.Cells.Item(9, 15).Value = 0.005363749481003509
$ws.Cells.Item(9, 16).Value = 0.005510546431623061
$ws.Cells.Item(9, 17).Value = 0.002528374212
$ws.Cells.Item(9, 18).Value = 0.022755367908
$ws.Cells.Item(9, 19).Value = [double]"2.169551246727999E-05"
$ws.Cells.Item(9, 20).Value = [double]"2.328223727078459E-05"

$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Efnb2"
$ws.Cells.Item(10, 3).Value = "Ephb1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.229894
$ws.Cells.Item(10, 8).Value = 0.6896819999999999
$ws.Cells.Item(10, 9).Value = 0.004044840748830231
$ws.Cells.Item(10, 10).Value = 0.004225032410066655
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.163866
$ws.Cells.Item(10, 14).Value = 0.327732
$ws.Cells.Item(10, 15).Value = 0.07991781891745053
$ws.Cells.Item(10, 16).Value = 0.05473669161449624
$ws.Cells.Item(10, 17).Value = 0.037671810204
$ws.Cells.Item(10, 18).Value = 0.226030861224
$ws.Cells.Item(10, 19).Value = 0.0003232548505149394
$ws.Cells.Item(10, 20).Value = 0.0002312642960910703

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Efnb2"
$ws.Cells.Item(11, 3).Value = "Ephb1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.6310036666666666
$ws.Cells.Item(11, 8).Value = 1.893011
$ws.Cells.Item(11, 9).Value = 0.01110211377241086
$ws.Cells.Item(11, 10).Value = 0.011596696488545
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.875567333333333
$ws.Cells.Item(11, 14).Value = 5.626702
$ws.Cells.Item(11, 15).Value = 0.9147184316015459
$ws.Cells.Item(11, 16).Value = 0.9397527619538806
$ws.Cells.Item(11, 17).Value = 1.183489864413555
$ws.Cells.Item(11, 18).Value = 10.651408779722
$ws.Cells.Item(11, 19).Value = 0.01015530809736158
$ws.Cells.Item(11, 20).Value = 0.01089802755465103

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Efnb2"
$ws.Cells.Item(12, 3).Value = "Ephb1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.6310036666666666
$ws.Cells.Item(12, 8).Value = 1.893011
$ws.Cells.Item(12, 9).Value = 0.01110211377241086
$ws.Cells.Item(12, 10).Value = 0.011596696488545
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.010998
$ws.Cells.Item(12, 14).Value = 0.032994
$ws.Cells.Item(12, 15).Value = 0.005363749481003509
$ws.Cells.Item(12, 16).Value = 0.005510546431623061
$ws.Cells.Item(12, 17).Value = 0.006939778326
$ws.Cells.Item(12, 18).Value = 0.06245800493400001
$ws.Cells.Item(12, 19).Value = [double]"5.954895698481064E-05"
$ws.Cells.Item(12, 20).Value = [double]"6.390413445356732E-05"

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Efnb2"
$ws.Cells.Item(13, 3).Value = "Ephb1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.6310036666666666
$ws.Cells.Item(13, 8).Value = 1.893011
$ws.Cells.Item(13, 9).Value = 0.01110211377241086
$ws.Cells.Item(13, 10).Value = 0.011596696488545
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.163866
$ws.Cells.Item(13, 14).Value = 0.327732
$ws.Cells.Item(13, 15).Value = 0.07991781891745053
$ws.Cells.Item(13, 16).Value = 0.05473669161449624
$ws.Cells.Item(13, 17).Value = 0.103400046842
$ws.Cells.Item(13, 18).Value = 0.6204002810520001
$ws.Cells.Item(13, 19).Value = 0.0008872567180644644
$ws.Cells.Item(13, 20).Value = 0.0006347647994403989

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Efnb2"
$ws.Cells.Item(14, 3).Value = "Ephb1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 7.2719705
$ws.Cells.Item(14, 8).Value = 14.543941
$ws.Cells.Item(14, 9).Value = 0.1279457602316344
$ws.Cells.Item(14, 10).Value = 0.08909703616318426
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.875567333333333
$ws.Cells.Item(14, 14).Value = 5.626702
$ws.Cells.Item(14, 15).Value = 0.9147184316015459
$ws.Cells.Item(14, 16).Value = 0.9397527619538806
$ws.Cells.Item(14, 17).Value = 13.63907031876367
$ws.Cells.Item(14, 18).Value = 81.834421912582
$ws.Cells.Item(14, 19).Value = 0.117034345129148
$ws.Cells.Item(14, 20).Value = 0.08372918581625718

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Efnb2"
$ws.Cells.Item(15, 3).Value = "Ephb1"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 7.2719705
$ws.Cells.Item(15, 8).Value = 14.543941
$ws.Cells.Item(15, 9).Value = 0.1279457602316344
$ws.Cells.Item(15, 10).Value = 0.08909703616318426
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.010998
$ws.Cells.Item(15, 14).Value = 0.032994
$ws.Cells.Item(15, 15).Value = 0.005363749481003509
$ws.Cells.Item(15, 16).Value = 0.005510546431623061
$ws.Cells.Item(15, 17).Value = 0.07997713155900001
$ws.Cells.Item(15, 18).Value = 0.4798627893540001
$ws.Cells.Item(15, 19).Value = 0.0006862690050390282
$ws.Cells.Item(15, 20).Value = 0.0004909733546972258

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Efnb2"
$ws.Cells.Item(16, 3).Value = "Ephb1"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 7.2719705
$ws.Cells.Item(16, 8).Value = 14.543941
$ws.Cells.Item(16, 9).Value = 0.1279457602316344
$ws.Cells.Item(16, 10).Value = 0.08909703616318426
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.163866
$ws.Cells.Item(16, 14).Value = 0.327732
$ws.Cells.Item(16, 15).Value = 0.07991781891745053
$ws.Cells.Item(16, 16).Value = 0.05473669161449624
$ws.Cells.Item(16, 17).Value = 1.191628717953
$ws.Cells.Item(16, 18).Value = 4.766514871812
$ws.Cells.Item(16, 19).Value = 0.0102251460974473
$ws.Cells.Item(16, 20).Value = 0.004876876992229836
